# Doing Updates for Financials
# Insert a new "most recent fiscal year" column (FY2018, period ending 2018-12-31)
# before the existing column D, shifting all the historical data one column to
# the right, then populate the new column with the FY2018 figures and apply a
# handful of restated prior-year figures that came in with this refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank column at D - this shifts the old D:K columns to E:L and
#    extends all the row "spans" / sheet dimension automatically.
$ws.Columns("D").Insert()

# 2) The newly inserted column D has no formatting of its own; clone the
#    number formats / fonts/ styles from column E (which now holds what used
#    to be column D, i.e. the same per-row style pattern) so the new FY2018
#    column looks like the rest of the table (date format on the header rows,
#    number format on the data rows).
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null

# 3) Populate the new column D with the FY2018 (period ending 2018-12-31)
#    figures for each of the three statements (Income Statement, Balance
#    Sheet, Cash Flow Statement).

# --- Income Statement (new FY2018 column) ---
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1242800
$ws.Range("D9").Value = 664500
$ws.Range("D10").Value = 578400
$ws.Range("D12").Value = 40800
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1114400
$ws.Range("D18").Value = 128400
$ws.Range("D20").Value = 8400
$ws.Range("D21").Value = 156800
$ws.Range("D22").Value = 5500
$ws.Range("D23").Value = 131300
$ws.Range("D24").Value = 26000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 105300
$ws.Range("D27").Value = 104700
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -8400
$ws.Range("D33").Value = 104700
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 104700

# --- Balance Sheet (new FY2018 column) ---
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 64000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 72100
$ws.Range("D44").Value = 338100
$ws.Range("D45").Value = 50800
$ws.Range("D46").Value = 524900
$ws.Range("D47").Value = 72200
$ws.Range("D48").Value = 88500
$ws.Range("D49").Value = 280500
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 86800
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1052900
$ws.Range("D57").Value = 42500
$ws.Range("D58").Value = 42900
$ws.Range("D59").Value = 217900
$ws.Range("D60").Value = 303200
$ws.Range("D61").Value = 7600
$ws.Range("D62").Value = 7800
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 328400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 413800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 724600
$ws.Range("D77").Value = 0

# --- Cash Flow Statement (new FY2018 column) ---
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 104700
$ws.Range("D83").Value = 19900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 92300
$ws.Range("D91").Value = -36800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -38500
$ws.Range("D96").Value = -3800
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -75100
$ws.Range("D101").Value = -400
$ws.Range("D102").Value = -21700

# 4) A handful of prior-year (now shifted) figures were also restated in this
#    refresh - fix those up explicitly.
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 0

$ws.Range("E45").Value = 22600
$ws.Range("F45").Value = 16600
$ws.Range("G45").Value = 16700
$ws.Range("I45").Value = 28900
$ws.Range("J45").Value = 30400
